$wb = $excel.ActiveWorkbook

function Set-CellValue {
    param($ws, [string]$addr, [double]$val)
    $ws.Range($addr).Value = $val
}

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$updates = @{
    "H33" = 306
    "I33" = 248.33333
    "K33" = 248.33333
    "M33" = -19.33332999999999
    "H43" = 1067
    "I43" = 599.5
    "K43" = 599.5
    "M43" = -530.5
    "H53" = 966.6111
    "I53" = 1101.8667
    "J53" = 290.33334
    "K53" = 1101.8667
    "L53" = 290.33334
    "M53" = -464.8667
    "N53" = -1564.33334
    "H64" = 200004740
    "I64" = 4999.6665
    "K64" = 4999.6665
    "M64" = -4751.6665
    "H67" = 200004740
    "I67" = 4999.6665
    "K67" = 4999.6665
    "M67" = -4141.6665
    "H76" = 3750
    "I76" = 3500
    "K76" = 3500
    "M76" = -3185
    "H79" = 3750
    "I79" = 3500
    "K79" = 3500
    "M79" = -2408
    "H97" = 2825
    "J97" = 2825
    "L97" = 8475
    "N97" = -9467
    "H111" = 2968.111
    "I111" = 4183.5
    "J111" = 1448.875
    "K111" = 12550.5
    "L111" = 4346.625
    "M111" = -9483.5
    "N111" = -10480.625
    "H113" = 30306194
    "I113" = 37039570
    "K113" = 37039570
    "M113" = -37036316
    "H132" = 6479.4546
    "I132" = 6952.1
    "J132" = 1753
    "K132" = 20856.3
    "L132" = 5259
    "M132" = -18326.3
    "N132" = -10319
    "H135" = 318
    "I135" = 345.66666
    "J135" = 110.5
    "K135" = 3110.99994
    "L135" = 994.5
    "M135" = -575.9999399999997
    "N135" = -6064.5
    "H137" = 2393.3462
    "I137" = 2392.9167
    "K137" = 7178.750100000001
    "M137" = -4628.750100000001
    "H138" = 1908.375
    "I138" = 964.1539
    "K138" = 2892.4617
    "M138" = 2247.5383
    "H141" = 2447.3125
    "I141" = 2373.6667
    "K141" = 7121.000100000001
    "M141" = -1941.000100000001
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$updates = @{
    "H32" = 1662067.6
    "I32" = 773858.4399999999
    "J32" = 15873416
    "K32" = 773858.4399999999
    "L32" = 15873416
    "M32" = -773571.4399999999
    "N32" = -15873990
    "H33" = 15666.667
    "I33" = 16000
    "K33" = 16000
    "M33" = -15671
    "H37" = 235
    "I37" = 235
    "J37" = 0
    "K37" = 235
    "L37" = 0
    "M37" = 38
    "H74" = 1795.9131
    "I74" = 1399.6666
    "K74" = 1399.6666
    "M74" = -525.6666
    "H77" = 1795.9131
    "I77" = 1399.6666
    "K77" = 6998.333000000001
    "M77" = -2630.333000000001
    "H95" = 36633.332
    "I95" = 20000
    "K95" = 20000
    "M95" = -17254
    "H110" = 815.8946999999999
    "I110" = 697.9167
    "K110" = 697.9167
    "M110" = 1347.0833
    "H132" = 4693.647
    "I132" = 6685.4287
    "J132" = 3299.4
    "K132" = 20056.2861
    "L132" = 9898.200000000001
    "M132" = -17526.2861
    "N132" = -14958.2
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

@("N37") | ForEach-Object { $ws.Range($_).ClearContents() }

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$updates = @{
    "H86" = 4990.7617
    "I86" = 4851.5
    "J86" = 5269.2856
    "K86" = 4851.5
    "L86" = 5269.2856
    "M86" = -3728.5
    "N86" = -7515.2856
    "H89" = 4990.7617
    "I89" = 4851.5
    "J89" = 5269.2856
    "K89" = 24257.5
    "L89" = 26346.428
    "M89" = -18641.5
    "N89" = -37578.428
    "H132" = 0
    "J132" = 0
    "L132" = 0
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

@("N132") | ForEach-Object { $ws.Range($_).ClearContents() }

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$updates = @{
    "H16" = 1538
    "I16" = 1538
    "K16" = 1538
    "M16" = -1251
    "H31" = 0
    "I31" = 0
    "J31" = 0
    "K31" = 0
    "L31" = 0
    "H34" = 0
    "I34" = 0
    "J34" = 0
    "K34" = 0
    "L34" = 0
    "H58" = 2252.8
    "I58" = 1632.5
    "K58" = 1632.5
    "M58" = -1429.5
    "H97" = 43598
    "J97" = 43598
    "L97" = 43598
    "N97" = -45580
    "H113" = 1538
    "I113" = 1538
    "K113" = 1538
    "M113" = 632
    "H122" = 4073.75
    "I122" = 3919
    "K122" = 11757
    "M122" = -9307
    "H125" = 50000
    "J125" = 50000
    "L125" = 50000
    "N125" = -54920
    "H132" = 3526.9487
    "I132" = 2938.4075
    "K132" = 8815.2225
    "M132" = -6285.2225
    "H134" = 3535.6287
    "I134" = 3863
    "K134" = 11589
    "M134" = -9054
    "H136" = 2252.8
    "I136" = 1632.5
    "K136" = 4897.5
    "M136" = -2347.5
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

@("M31","N31","M34","N34") | ForEach-Object { $ws.Range($_).ClearContents() }

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$updates = @{
    "H109" = 3242.1667
    "I109" = 1044.875
    "J109" = 5000
    "K109" = 3134.625
    "L109" = 15000
    "M109" = -2094.625
    "N109" = -17080
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$updates = @{
    "H70" = 5619.091
    "I70" = 4776.625
    "J70" = 7865.6665
    "K70" = 4776.625
    "L70" = 7865.6665
    "M70" = -4506.625
    "N70" = -8405.666499999999
    "H73" = 5619.091
    "I73" = 4776.625
    "J73" = 7865.6665
    "K73" = 4776.625
    "L73" = 7865.6665
    "M73" = -3840.625
    "N73" = -9737.666499999999
    "H102" = 1306.9773
    "I102" = 1000.45
    "K102" = 1000.45
    "M102" = 621.55
    "H116" = 49999
    "J116" = 49999
    "L116" = 49999
    "N116" = -59177
    "H122" = 3150.6428
    "I122" = 2291.1428
    "J122" = 4010.1428
    "K122" = 6873.428400000001
    "L122" = 12030.4284
    "M122" = -4423.428400000001
    "N122" = -16930.4284
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$updates = @{
    "H132" = 4373.3335
    "I132" = 5407.778
    "K132" = 16223.334
    "M132" = -13693.334
    "H136" = 8584
    "I136" = 10749.5
    "J136" = 7501.25
    "K136" = 32248.5
    "L136" = 22503.75
    "M136" = -29698.5
    "N136" = -27603.75
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$updates = @{
    "H136" = 3910.2354
    "J136" = 3024.125
    "L136" = 9072.375
    "N136" = -14172.375
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Host "Applied all Leve profit updates across sheets."